$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from D1 to E1 so the new header matches the existing ones
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("E1").Value = "Colocação"

$ws.Range("E2").Value = "1º"
$ws.Range("E3").Value = "2º"
$ws.Range("E4").Value = "3º"
$ws.Range("E5").Value = "4º"
$ws.Range("E6").Value = "5º"
$ws.Range("E7").Value = "6º"
$ws.Range("E8").Value = "23º"
